$d = $word.ActiveDocument

# 1. Replace the court-venue-name placeholder with the external-short-name
#    placeholder in the heading ("In the County Court at <<...>>").
$d.Content.Find.Execute(
    ".venue_name>><<else>> Online Civil Claims<<es_>>",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ".external_short_name>><<else>>Online Civil Claims<<es_>>", 2) | Out-Null

# 2. The heading paragraph now wraps onto a second line because the new
#    placeholder text is longer, so Word (re)applies 1.15-line spacing to it.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*In the County Court at*") {
        $p.Range.ParagraphFormat.LineSpacingRule = 5
        $p.Range.ParagraphFormat.LineSpacing = 13.8
    }
}
